$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "Title" parameter label used for the STITLE rows should read "Study Title"
# (CreateTsXPT.R / write.xport2.R now require the longer, more descriptive label).
$ws.Range("F6").Value = "Study Title"
$ws.Range("F10").Value = "Study Title"
$ws.Range("F14").Value = "Study Title"

# Reflect the current selection/active cell as recorded in the saved workbook.
$ws.Activate()
$ws.Range("F15").Select()
